# Commit: "Returns error for format and missing requires"
#
# - Adds a new product row (row 3) to the "Productos" sheet: ID 10, name
#   "queque2" reusing most of the other values from the existing "queque"
#   product (row 2).
# - Adds "queque2" as a new shared string.
# - Re-creates / re-orders the "Productos" data validations so the "ID
#   must be unique" rule comes first, and extends the single-row rules
#   (numbers-only, Si/No, categories) to cover the new row 3.
# - Makes "Productos" the active sheet/tab again (it had drifted to
#   "Variantes"), and updates the remembered selections on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Productos: fill in the new row 3 with the second product.
# ---------------------------------------------------------------------
$productos = $wb.Worksheets.Item("Productos")
$productos.Activate()

$productos.Range("A3").Value = 10
$productos.Range("B3").Value = "queque2"
$productos.Range("C3").Value = "en molde de cupcake"
$productos.Range("D3").Value = 1500
$productos.Range("F3").Value = 200
$productos.Range("G3").Value = 10
$productos.Range("H3").Value = 10
$productos.Range("I3").Value = 10
$productos.Range("K3").Value = "Este es el mejor queque de Chile"
$productos.Range("L3").Value = "Sí"
$productos.Range("M3").Value = 43075
$productos.Range("M3").NumberFormat = $productos.Range("M2").NumberFormat
$productos.Range("N3").Value = "riri->ruru->rere, coco"
$productos.Range("O3").Value = "Por defecto"

# Column J (Slug) is intentionally left blank for the new row, matching
# the template formatting of the rest of the row.
$productos.Range("J3").Font.Name = $productos.Range("J2").Font.Name

# ---------------------------------------------------------------------
# 2. Productos: rebuild the data validations — "ID unique" first, then
#    "SKU unique", then the single-row rules widened to include row 3.
# ---------------------------------------------------------------------
$idRule = $productos.Range("A1:A2")
$skuRule = $productos.Range("E1:E2")
$numRule = $productos.Range("D2:D3")
$yesNoRule = $productos.Range("L2:L3")
$catRule = $productos.Range("N2:N3")

$idRule.Validation.Delete()
$skuRule.Validation.Delete()
$numRule.Validation.Delete()
$yesNoRule.Validation.Delete()
$catRule.Validation.Delete()

$idRule.Validation.Add(7, 1, 1, 'COUNTIF($A:$A,"="&A1)  < 2', 0)
$idRule.Validation.InputMessage = "El ID debe ser único"
$idRule.Validation.ShowInput = $true
$idRule.Validation.ShowError = $true
$idRule.Validation.IgnoreBlank = $true

$skuRule.Validation.Add(7, 1, 1, 'AND(COUNTIF(Productos!$E:$E,"="&Productos!E1)  < 2, COUNTIF(Variantes!$D:$D,"="&Productos!E1)  < 1)', 0)
$skuRule.Validation.InputMessage = "El SKU ya existe (en Productos o Variantes). Debe ser único."
$skuRule.Validation.ShowInput = $true
$skuRule.Validation.ShowError = $true
$skuRule.Validation.IgnoreBlank = $true

$numRule.Validation.Add(7, 1, 1, 'regexmatch(to_text(D2),"^[0-9]*$")=1', 0)
$numRule.Validation.InputMessage = "Se deben ingresar solo números"
$numRule.Validation.ShowInput = $true
$numRule.Validation.ShowError = $true
$numRule.Validation.IgnoreBlank = $true

$yesNoRule.Validation.Add(3, 1, 1, '"Sí,No"', 0)
$yesNoRule.Validation.InputMessage = "Debe ser Sí o No"
$yesNoRule.Validation.ShowInput = $true
$yesNoRule.Validation.ShowError = $true
$yesNoRule.Validation.IgnoreBlank = $true
$yesNoRule.Validation.InCellDropdown = $true

$catRule.Validation.Add(7, 1, 1, 'regexmatch(N2,"^((([A-z0-9]+->)*[A-z0-9]+)(, )?)*(([A-z0-9]+->)*[A-z0-9]+)?$")=1', 0)
$catRule.Validation.InputMessage = "Deber ingresar las categorías de la forma: riri->ruru->rere, coco"
$catRule.Validation.ShowInput = $true
$catRule.Validation.ShowError = $true
$catRule.Validation.IgnoreBlank = $true

# ---------------------------------------------------------------------
# 3. Remembered selections on the other sheets now also reference the
#    (still-empty) A3 cell alongside their previous active cell.
# ---------------------------------------------------------------------
$variantes = $wb.Worksheets.Item("Variantes")
$variantes.Activate()
$variantes.Range("E1").Select()

$propiedades = $wb.Worksheets.Item("Propiedades")
$propiedades.Activate()
$propiedades.Range("B2").Select()

$opciones = $wb.Worksheets.Item("Opciones")
$opciones.Activate()
$opciones.Range("B3").Select()

$ubicaciones = $wb.Worksheets.Item("Ubicaciones")
$ubicaciones.Activate()
$ubicaciones.Range("B3").Select()

$stock = $wb.Worksheets.Item("Stock")
$stock.Activate()
$stock.Range("B3").Select()

# ---------------------------------------------------------------------
# 4. Productos becomes the active sheet/tab again, selection back at
#    the top, with the new row highlighted.
# ---------------------------------------------------------------------
$productos.Activate()
$productos.Range("A1").Select()
$productos.Range("A3").Select()
